# Insert a new weekly price record as row 35 in the "Fruta, Vega Modelo de
# Temuco - Membrillo" sheet. All existing rows from 35 downward shift down
# by one (row 35 -> 36, ..., row 127 -> 128), keeping their original data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 35..127 down by one row, creating a blank row 35.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Cells.Item(35, 1).Value  = 10
$ws.Cells.Item(35, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value  = "La Araucanía"
$ws.Cells.Item(35, 4).Value  = 44623
$ws.Cells.Item(35, 5).Value  = 9
$ws.Cells.Item(35, 6).Value  = "Fruta"
$ws.Cells.Item(35, 7).Value  = 100104
$ws.Cells.Item(35, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(35, 9).Value  = 100104003
$ws.Cells.Item(35, 10).Value = "Membrillo"
$ws.Cells.Item(35, 11).Value = "Champion"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 90
$ws.Cells.Item(35, 14).Value = 17000
$ws.Cells.Item(35, 15).Value = 18000
$ws.Cells.Item(35, 16).Value = 17389
$ws.Cells.Item(35, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(35, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(35, 19).Value = 966
$ws.Cells.Item(35, 20).Value = 18
